$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 470.25
$ws.Range("I33").Value = 354.10257
$ws.Range("K33").Value = 354.10257
$ws.Range("M33").Value = -125.10257

$ws.Range("H40").Value = 14224.25
$ws.Range("I40").Value = 18299
$ws.Range("K40").Value = 18299
$ws.Range("M40").Value = -18124

$ws.Range("H70").Value = 53031348
$ws.Range("J70").Value = 41667790
$ws.Range("L70").Value = 125003370
$ws.Range("N70").Value = -125003910

$ws.Range("H73").Value = 53031348
$ws.Range("J73").Value = 41667790
$ws.Range("L73").Value = 125003370
$ws.Range("N73").Value = -125005242

$ws.Range("H81").Value = 43000
$ws.Range("J81").Value = 43000
$ws.Range("L81").Value = 43000
$ws.Range("N81").Value = -44996

$ws.Range("H84").Value = 43000
$ws.Range("J84").Value = 43000
$ws.Range("L84").Value = 129000
$ws.Range("N84").Value = -138984

$ws.Range("H112").Value = 9923.913
$ws.Range("J112").Value = 10234.091
$ws.Range("L112").Value = 30702.273
$ws.Range("N112").Value = -32918.273

$ws.Range("H113").Value = 33342100
$ws.Range("I113").Value = 6285.5713
$ws.Range("K113").Value = 6285.5713
$ws.Range("M113").Value = -3031.5713

$ws.Range("H116").Value = 13161356
$ws.Range("I116").Value = 19233266
$ws.Range("J116").Value = 5549.1665
$ws.Range("K116").Value = 19233266
$ws.Range("L116").Value = 5549.1665
$ws.Range("M116").Value = -19229824
$ws.Range("N116").Value = -12433.1665

$ws.Range("H129").Value = 1358.5834
$ws.Range("J129").Value = 2318.3333
$ws.Range("L129").Value = 6954.999899999999
$ws.Range("N129").Value = -16954.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1545240.6
$ws.Range("I32").Value = 1668730
$ws.Range("K32").Value = 1668730
$ws.Range("M32").Value = -1668443

$ws.Range("H61").Value = 6943.122
$ws.Range("I61").Value = 3351
$ws.Range("J61").Value = 10714.85
$ws.Range("K61").Value = 3351
$ws.Range("L61").Value = 10714.85
$ws.Range("M61").Value = -3139
$ws.Range("N61").Value = -11138.85

$ws.Range("H74").Value = 50540.316
$ws.Range("I74").Value = 70339.266
$ws.Range("J74").Value = 8114
$ws.Range("K74").Value = 70339.266
$ws.Range("L74").Value = 8114
$ws.Range("M74").Value = -69465.266
$ws.Range("N74").Value = -9862

$ws.Range("H77").Value = 50540.316
$ws.Range("I77").Value = 70339.266
$ws.Range("J77").Value = 8114
$ws.Range("K77").Value = 351696.33
$ws.Range("L77").Value = 40570
$ws.Range("M77").Value = -347328.33
$ws.Range("N77").Value = -49306

$ws.Range("H136").Value = 6943.122
$ws.Range("I136").Value = 3351
$ws.Range("J136").Value = 10714.85
$ws.Range("K136").Value = 10053
$ws.Range("L136").Value = 32144.55
$ws.Range("M136").Value = -7503
$ws.Range("N136").Value = -37244.55

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 15153756
$ws.Range("I20").Value = 23810986
$ws.Range("K20").Value = 23810986
$ws.Range("M20").Value = -23810739

$ws.Range("H94").Value = 41669490
$ws.Range("I94").Value = 66668508
$ws.Range("K94").Value = 66668508
$ws.Range("M94").Value = -66668057

$ws.Range("H99").Value = 5685749
$ws.Range("I99").Value = 2999.6
$ws.Range("K99").Value = 2999.6
$ws.Range("M99").Value = -1501.6

$ws.Range("H107").Value = 56253028
$ws.Range("I107").Value = 86540550
$ws.Range("J107").Value = 4771
$ws.Range("K107").Value = 86540550
$ws.Range("L107").Value = 4771
$ws.Range("M107").Value = -86538630
$ws.Range("N107").Value = -8611

$ws.Range("H134").Value = 4765.718
$ws.Range("J134").Value = 12262
$ws.Range("L134").Value = 36786
$ws.Range("N134").Value = -41856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("N44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("N44").ClearContents()

$ws.Range("H45").Value = 22415
$ws.Range("J45").Value = 22415
$ws.Range("L45").Value = 22415
$ws.Range("N45").Value = -23601

$ws.Range("H62").Value = 1665.6666
$ws.Range("I62").Value = 1499.5
$ws.Range("J62").Value = 1998
$ws.Range("K62").Value = 1499.5
$ws.Range("L62").Value = 1998
$ws.Range("M62").Value = -875.5
$ws.Range("N62").Value = -3246

$ws.Range("H65").Value = 1665.6666
$ws.Range("I65").Value = 1499.5
$ws.Range("J65").Value = 1998
$ws.Range("K65").Value = 7497.5
$ws.Range("L65").Value = 9990
$ws.Range("M65").Value = -4377.5
$ws.Range("N65").Value = -16230

$ws.Range("H99").Value = 7699.9414
$ws.Range("J99").Value = 8460.134
$ws.Range("L99").Value = 8460.134
$ws.Range("N99").Value = -11456.134

$ws.Range("H107").Value = 1028.68
$ws.Range("J107").Value = 2700.6667
$ws.Range("L107").Value = 2700.6667
$ws.Range("N107").Value = -6540.6667

$ws.Range("H126").Value = 7699.9414
$ws.Range("J126").Value = 8460.134
$ws.Range("L126").Value = 25380.402
$ws.Range("N126").Value = -30320.402

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2366.28
$ws.Range("I68").Value = 1832.3334
$ws.Range("J68").Value = 2534.8948
$ws.Range("K68").Value = 5497.0002
$ws.Range("L68").Value = 7604.6844
$ws.Range("M68").Value = -4686.0002
$ws.Range("N68").Value = -9226.6844

$ws.Range("H71").Value = 2366.28
$ws.Range("I71").Value = 1832.3334
$ws.Range("J71").Value = 2534.8948
$ws.Range("K71").Value = 16491.0006
$ws.Range("L71").Value = 22814.0532
$ws.Range("M71").Value = -12435.0006
$ws.Range("N71").Value = -30926.0532

$ws.Range("H131").Value = 1871.66
$ws.Range("I131").Value = 1022
$ws.Range("J131").Value = 2111.3076
$ws.Range("K131").Value = 3066
$ws.Range("L131").Value = 6333.9228
$ws.Range("M131").Value = 1974
$ws.Range("N131").Value = -16413.9228

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 24417608
$ws.Range("I122").Value = 43524052
$ws.Range("J122").Value = 3815.7778
$ws.Range("K122").Value = 130572156
$ws.Range("L122").Value = 11447.3334
$ws.Range("M122").Value = -130569706
$ws.Range("N122").Value = -16347.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7063.4814
$ws.Range("I7").Value = 4213.125
$ws.Range("K7").Value = 4213.125
$ws.Range("M7").Value = -4101.125

$ws.Range("H22").Value = 1330.2632
$ws.Range("I22").Value = 977.7857
$ws.Range("K22").Value = 977.7857
$ws.Range("M22").Value = -682.7857

$ws.Range("H27").Value = 1330.2632
$ws.Range("I27").Value = 977.7857
$ws.Range("K27").Value = 977.7857
$ws.Range("M27").Value = -870.7857

$ws.Range("H36").Value = 59999.668
$ws.Range("I36").Value = 59999
$ws.Range("K36").Value = 59999
$ws.Range("M36").Value = -59437

$ws.Range("H40").Value = 125005250
$ws.Range("I40").Value = 166669660
$ws.Range("J40").Value = 12000
$ws.Range("K40").Value = 166669660
$ws.Range("L40").Value = 12000
$ws.Range("M40").Value = -166669524
$ws.Range("N40").Value = -12272

$ws.Range("H55").Value = 1284.65
$ws.Range("I55").Value = 1440.2307
$ws.Range("J55").Value = 995.7143
$ws.Range("K55").Value = 1440.2307
$ws.Range("L55").Value = 995.7143
$ws.Range("M55").Value = -1267.2307
$ws.Range("N55").Value = -1341.7143

$ws.Range("H124").Value = 51528
$ws.Range("J124").Value = 51528
$ws.Range("L124").Value = 51528
$ws.Range("N124").Value = -61348

$ws.Range("H126").Value = 7063.4814
$ws.Range("I126").Value = 4213.125
$ws.Range("K126").Value = 12639.375
$ws.Range("M126").Value = -10169.375

$ws.Range("H128").Value = 67905.8
$ws.Range("J128").Value = 67905.8
$ws.Range("L128").Value = 67905.8
$ws.Range("N128").Value = -77865.8

$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("N129").Value = 0
$ws.Range("N129").ClearContents()

$ws.Range("H134").Value = 83990
$ws.Range("J134").Value = 83990
$ws.Range("L134").Value = 83990
$ws.Range("N134").Value = -94130

$ws.Range("H136").Value = 9739.098
$ws.Range("I136").Value = 7061.8
$ws.Range("J136").Value = 12288.904
$ws.Range("K136").Value = 21185.4
$ws.Range("L136").Value = 36866.712
$ws.Range("M136").Value = -18635.4
$ws.Range("N136").Value = -41966.712

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 16159783
$ws.Range("I81").Value = 910653.2
$ws.Range("K81").Value = 1821306.4
$ws.Range("M81").Value = -1820245.4

$ws.Range("H84").Value = 16159783
$ws.Range("I84").Value = 910653.2
$ws.Range("K84").Value = 9106532
$ws.Range("M84").Value = -9101228

$ws.Range("H122").Value = 7103579
$ws.Range("I122").Value = 10083247
$ws.Range("K122").Value = 30249741
$ws.Range("M122").Value = -30247291

$ws.Range("H126").Value = 4299.95
$ws.Range("I126").Value = 1550.2
$ws.Range("J126").Value = 7049.7
$ws.Range("K126").Value = 4650.6
$ws.Range("L126").Value = 21149.1
$ws.Range("M126").Value = -2180.6
$ws.Range("N126").Value = -26089.1

$ws.Range("H132").Value = 9685.370999999999
$ws.Range("J132").Value = 17507.264
$ws.Range("L132").Value = 52521.792
$ws.Range("N132").Value = -57581.792

$ws.Range("H136").Value = 28502.926
$ws.Range("I136").Value = 1463.7667
$ws.Range("K136").Value = 4391.300099999999
$ws.Range("M136").Value = -1841.300099999999
